# Feat: Add more test case
# Insert 5 new rows before row 23 (pushes existing cleanup rows 23-27 down to 28-32),
# fill the newly freed rows 23-27 with the new "test5user" test case flow, update the
# DELAY value of the last row, and refresh the sheet selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: insert 5 blank rows starting at row 23. This shifts the previous
# rows 23-27 (DELETE testuser / test2user / test3user / testwronguser, DELAY 5)
# down to rows 28-32.
$ws.Range("A23:A27").EntireRow.Insert()

# Populate the new payload/endpoint strings first (matches the order the new
# shared-string table entries were authored in).
$ws.Range("C23").Value = '{"UserName": "test5user", "Password": "Test1234!", "RoleId": "Administrator"}'
$ws.Range("C24").Value = '{"PasswordChangeRequired":true}'
$ws.Range("C25").Value = '{"UserName":"test5user","Password":"Test1234!"}'
$ws.Range("B24").Value = '/redfish/v1/AccountService/Accounts/${test5user.id}'
$ws.Range("C26").Value = '{"AccountLockoutThreshold":4}'

# Row 23: create the new "test5user" account.
$ws.Range("A23").Value = "POST"
$ws.Range("B23").Value = "/redfish/v1/AccountService/Accounts"

# Row 24: force the new user to require a password change.
$ws.Range("A24").Value = "PATCH"

# Row 25: log in as the new user (styled the same as the other login row, with
# wrapped text and a taller row height).
$ws.Range("A25").Value = "POST"
$ws.Range("B25").Value = "/redfish/v1/SessionService/Sessions"
$ws.Range("C25").WrapText = $true
$ws.Rows.Item(25).RowHeight = 16

# Row 26: set an account lockout threshold.
$ws.Range("A26").Value = "PATCH"
$ws.Range("B26").Value = "/redfish/v1/AccountService"

# Row 27: clean up the new test5user account.
$ws.Range("A27").Value = "DELETE"
$ws.Range("B27").Value = '/redfish/v1/AccountService/Accounts/${test5user.id}'

# The pre-existing DELAY row (now row 32) changes its wait value from 5 to 3.
$ws.Range("B32").Value = 3

# Update the active selection to reflect where the editor cursor ended up.
$ws.Range("C28").Select()
